$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Capture values that need to move before we clear anything
$analisarCodigo = $ws.Range("D4").Value2
$darIdeias      = $ws.Range("E4").Value2
$jogarOJogo     = $ws.Range("E5").Value2
$pesquisarJogo  = $ws.Range("E6").Value2

# Clear the left (user story) column B4:B6 - wrongly merged content
$ws.Range("B4").ClearContents()
$ws.Range("B5").ClearContents()
$ws.Range("B6").ClearContents()

# Clear D4 (value relocates to E4) and E5:E6 (values relocate to column G)
$ws.Range("D4").ClearContents()
$ws.Range("E5").ClearContents()
$ws.Range("E6").ClearContents()

# Move "Analisar o código dado" from D4 into E4
$ws.Range("E4").Value = $analisarCodigo

# Append the relocated "Doing" items to the bottom of the "Done" column (G)
$ws.Range("G9").Value  = "Fazer pdf dos 3 User Stories mais votados pela equipa e submeter no moodle"
$ws.Range("G10").Value = $pesquisarJogo
$ws.Range("G11").Value = $jogarOJogo
$ws.Range("G12").Value = $darIdeias

# Widen column G to fit the newly added content
$ws.Columns.Item(7).ColumnWidth = 68

# Update the active cell selection to reflect where editing finished
$ws.Range("G15").Select()
